# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''97.002.91'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -1.78%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''3.315.47'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -4.78%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  -0.14%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''246.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -5.79%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''650.16'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -4.18%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''1.36'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  -12.83%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.414'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -10.98%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.998'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -0.08%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''0.984'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -10.74%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''3.313.41'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -4.81%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.204'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -8.16%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''39.76'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -6.95%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''96.698.39'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -2.28%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''5.97'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -4.71%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''0.0000250'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -9.13%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''3.936.29'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -4.91%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''8.63'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +5.57%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''3.300.11'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -5.18%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''16.66'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -5.38%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''0.494'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +9.26%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''10.43'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -3.21%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''490.69'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -8.60%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''3.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -10.61%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''0.0000197'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -11.04%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''6.25'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -0.52%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''92.31'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -10.13%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''11.96'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -7.95%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''3.481.69'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  -5.09%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = '''Hedera'
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = '''0.142'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -5.74%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = '''Dai'
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = '''https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = '''1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +0.06%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''10.75'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -5.79%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''0.186'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -7.65%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''2.45'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  +9.39%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -0.02%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''0.542'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -7.73%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''27.94'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -9.54%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''1.46'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +2.41%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''7.46'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -6.71%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = '''  -0.06%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.149'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -7.75%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''501.37'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -7.57%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''24.51'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -1.04%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''3.70'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -1.69%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.819'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -5.89%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = '''VeChain'
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = '''0.0405'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -8.08%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = '''Cosmos'
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = '''https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = '''8.37'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +1.60%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''5.41'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +2.13%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''1.62'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  +1.17%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''52.44'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +1.98%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''3.10'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -11.86%  '
$ws.Range("E51").Style = "Normal"
